$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C5").Value = -45
$ws.Range("C6").Value = 125
$ws.Range("C7").Value = 40
$ws.Range("C8").Value = 45
$ws.Range("C12").Value = 100
$ws.Range("C13").Value = -3.5
$ws.Range("C14").Value = -92

$ws.Range("E4").Formula = "=C12+C2"
$ws.Range("F4").Formula = "=C13+C3"
$ws.Range("G4").Formula = "=C14+C4"
$ws.Range("J4").Formula = "=-E4*SIN(RADIANS(C5)) + G4*COS(RADIANS(C5))"

$excel.Calculate()
